$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.576.49'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '2.304.35'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.63'
$ws.Range('E5').Value = '  -1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.50'
$ws.Range('E6').Value = '  -4.54%  '
$ws.Range('E7').Value = '  -3.60%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.503'
$ws.Range('E9').Value = '  -3.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.74'
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').Value = '2.660.82'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.65'
$ws.Range('E15').Value = '  +4.05%  '
$ws.Range('D16').Value = '2.306.63'
$ws.Range('E16').Value = '  -11.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.804'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '42.520.89'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('D19').Value = '0.0₃0906'
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.08'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.46'
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.85'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.97'
$ws.Range('E23').Value = '  -2.21%  '
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.95'
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('E28').Value = '  +9.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.73'
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '163.89'
$ws.Range('E30').Value = '  +1.43%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.15'
$ws.Range('E31').Value = '  -3.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.01'
$ws.Range('E33').Value = '  -3.93%  '
$ws.Range('E34').Value = '  -4.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0703'
$ws.Range('E35').Value = '  -4.21%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.74'
$ws.Range('E36').Value = '  -7.88%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.42'
$ws.Range('E37').Value = '  -2.25%  '
$ws.Range('E38').Value = '  -3.73%  '
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('E40').Value = '  -5.58%  '
$ws.Range('E41').Value = '  -3.10%  '
$ws.Range('E42').Value = '  -3.51%  '
$ws.Range('D43').Value = '1.965.39'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.50'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.23'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('E47').Value = '  -5.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.65'
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('E49').Value = '  -3.44%  '
$ws.Range('D50').Value = '2.527.17'
$ws.Range('E50').Value = '  -0.06%  '
$ws.Range('E51').Value = '  +0.32%  '
